# Update the "Förändrad" (Changed) date column C for all data rows (2-27)
# from serial date 45319 (2024-01-28) to serial date 45321 (2024-01-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45319) {
        $cell.Value2 = 45321
    }
}
